$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1196.5333
$ws.Range("J41").Value = 1324.75
$ws.Range("L41").Value = 1324.75
$ws.Range("N41").Value = -2204.75
$ws.Range("H64").Value = 5049.9
$ws.Range("J64").Value = 7833.3335
$ws.Range("L64").Value = 7833.3335
$ws.Range("N64").Value = -8329.333500000001
$ws.Range("H67").Value = 5049.9
$ws.Range("J67").Value = 7833.3335
$ws.Range("L67").Value = 7833.3335
$ws.Range("N67").Value = -9549.333500000001
$ws.Range("H129").Value = 27778672
$ws.Range("I129").Value = 29412534
$ws.Range("K129").Value = 88237602
$ws.Range("M129").Value = -88232602
$ws.Range("H131").Value = 4988.846
$ws.Range("I131").Value = 4137
$ws.Range("J131").Value = 7828.3335
$ws.Range("K131").Value = 12411
$ws.Range("L131").Value = 23485.0005
$ws.Range("M131").Value = -7371
$ws.Range("N131").Value = -33565.00049999999
$ws.Range("H138").Value = 3342.551
$ws.Range("I138").Value = 1967.4117
$ws.Range("J138").Value = 4073.0938
$ws.Range("K138").Value = 5902.2351
$ws.Range("L138").Value = 12219.2814
$ws.Range("M138").Value = -762.2350999999999
$ws.Range("N138").Value = -22499.2814

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 56949
$ws.Range("I80").Value = 33999
$ws.Range("J80").Value = 79899
$ws.Range("K80").Value = 33999
$ws.Range("L80").Value = 79899
$ws.Range("M80").Value = -33001
$ws.Range("N80").Value = -81895
$ws.Range("H83").Value = 56949
$ws.Range("I83").Value = 33999
$ws.Range("J83").Value = 79899
$ws.Range("K83").Value = 101997
$ws.Range("L83").Value = 239697
$ws.Range("M83").Value = -97005
$ws.Range("N83").Value = -249681
$ws.Range("H110").Value = 8676.913
$ws.Range("I110").Value = 2813.5789
$ws.Range("K110").Value = 2813.5789
$ws.Range("M110").Value = -768.5789
$ws.Range("H122").Value = 9664495
$ws.Range("J122").Value = 5498.5
$ws.Range("L122").Value = 16495.5
$ws.Range("N122").Value = -21395.5
$ws.Range("H132").Value = 1604.6046
$ws.Range("I132").Value = 942.4375
$ws.Range("J132").Value = 3530.9092
$ws.Range("K132").Value = 2827.3125
$ws.Range("L132").Value = 10592.7276
$ws.Range("M132").Value = -297.3125
$ws.Range("N132").Value = -15652.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 30000
$ws.Range("J32").Value = 30000
$ws.Range("L32").Value = 30000
$ws.Range("N32").Value = -30768
$ws.Range("H105").Value = 2644.4285
$ws.Range("J105").Value = 1375
$ws.Range("L105").Value = 1375
$ws.Range("N105").Value = -4869

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 5424.75
$ws.Range("J21").Value = 5424.75
$ws.Range("L21").Value = 5424.75
$ws.Range("N21").Value = -5894.75
$ws.Range("H23").Value = 30000
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = ""
$ws.Range("H27").Value = 30000
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").Value = ""
$ws.Range("H31").Value = 41916.824
$ws.Range("I31").Value = 1958.7778
$ws.Range("J31").Value = 67604.14
$ws.Range("K31").Value = 1958.7778
$ws.Range("L31").Value = 67604.14
$ws.Range("M31").Value = -1663.7778
$ws.Range("N31").Value = -68194.14
$ws.Range("H34").Value = 41916.824
$ws.Range("I34").Value = 1958.7778
$ws.Range("J34").Value = 67604.14
$ws.Range("K34").Value = 1958.7778
$ws.Range("L34").Value = 67604.14
$ws.Range("M34").Value = -1756.7778
$ws.Range("N34").Value = -68008.14
$ws.Range("H99").Value = 3427.6428
$ws.Range("I99").Value = 2624.625
$ws.Range("J99").Value = 4498.3335
$ws.Range("K99").Value = 2624.625
$ws.Range("L99").Value = 4498.3335
$ws.Range("M99").Value = -1126.625
$ws.Range("N99").Value = -7494.3335
$ws.Range("H122").Value = 2192.9688
$ws.Range("I122").Value = 1711.7084
$ws.Range("J122").Value = 3636.75
$ws.Range("K122").Value = 5135.1252
$ws.Range("L122").Value = 10910.25
$ws.Range("M122").Value = -2685.1252
$ws.Range("N122").Value = -15810.25
$ws.Range("H126").Value = 3427.6428
$ws.Range("I126").Value = 2624.625
$ws.Range("J126").Value = 4498.3335
$ws.Range("K126").Value = 7873.875
$ws.Range("L126").Value = 13495.0005
$ws.Range("M126").Value = -5403.875
$ws.Range("N126").Value = -18435.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4839164
$ws.Range("I4").Value = 5293629
$ws.Range("K4").Value = 15880887
$ws.Range("M4").Value = -15880775
$ws.Range("H5").Value = 1076.8966
$ws.Range("I5").Value = 714.2353000000001
$ws.Range("J5").Value = 1590.6666
$ws.Range("K5").Value = 2142.7059
$ws.Range("L5").Value = 4771.9998
$ws.Range("M5").Value = -2030.7059
$ws.Range("N5").Value = -4995.9998
$ws.Range("H47").Value = 849.75
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = ""
$ws.Range("H56").Value = 20838824
$ws.Range("I56").Value = 20838824
$ws.Range("K56").Value = 20838824
$ws.Range("M56").Value = -20838294
$ws.Range("H96").Value = 14793.3
$ws.Range("I96").Value = 12933
$ws.Range("J96").Value = 15000
$ws.Range("K96").Value = 38799
$ws.Range("L96").Value = 45000
$ws.Range("M96").Value = -36740
$ws.Range("N96").Value = -49118
$ws.Range("H105").Value = 4837
$ws.Range("J105").Value = 4837
$ws.Range("L105").Value = 14511
$ws.Range("N105").Value = -19753
$ws.Range("H122").Value = 925.1667
$ws.Range("J122").Value = 934.4
$ws.Range("L122").Value = 8409.6
$ws.Range("N122").Value = -13309.6
$ws.Range("H132").Value = 1866.4706
$ws.Range("I132").Value = 1153.5555
$ws.Range("J132").Value = 2123.12
$ws.Range("K132").Value = 10381.9995
$ws.Range("L132").Value = 19108.08
$ws.Range("M132").Value = -7851.9995
$ws.Range("N132").Value = -24168.08
$ws.Range("H135").Value = 1076.8966
$ws.Range("I135").Value = 714.2353000000001
$ws.Range("J135").Value = 1590.6666
$ws.Range("K135").Value = 6428.117700000001
$ws.Range("L135").Value = 14315.9994
$ws.Range("M135").Value = -3893.117700000001
$ws.Range("N135").Value = -19385.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").Value = ""
$ws.Range("H122").Value = 361179.8
$ws.Range("J122").Value = 7580.909
$ws.Range("L122").Value = 22742.727
$ws.Range("N122").Value = -27642.727
$ws.Range("H126").Value = 3698.625
$ws.Range("I126").Value = 3349
$ws.Range("K126").Value = 10047
$ws.Range("M126").Value = -7577

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5189.8667
$ws.Range("I46").Value = 3937
$ws.Range("J46").Value = 6621.7144
$ws.Range("K46").Value = 3937
$ws.Range("L46").Value = 6621.7144
$ws.Range("M46").Value = -3749
$ws.Range("N46").Value = -6997.7144
$ws.Range("H82").Value = 1290.2142
$ws.Range("I82").Value = 921.125
$ws.Range("J82").Value = 1782.3334
$ws.Range("K82").Value = 921.125
$ws.Range("L82").Value = 1782.3334
$ws.Range("M82").Value = -560.125
$ws.Range("N82").Value = -2504.3334
$ws.Range("H85").Value = 1290.2142
$ws.Range("I85").Value = 921.125
$ws.Range("J85").Value = 1782.3334
$ws.Range("K85").Value = 921.125
$ws.Range("L85").Value = 1782.3334
$ws.Range("M85").Value = 326.875
$ws.Range("N85").Value = -4278.3334
$ws.Range("H132").Value = 6942.8667
$ws.Range("I132").Value = 7658.923
$ws.Range("J132").Value = 5963
$ws.Range("K132").Value = 22976.769
$ws.Range("L132").Value = 17889
$ws.Range("M132").Value = -20446.769
$ws.Range("N132").Value = -22949
$ws.Range("H136").Value = 72199.47
$ws.Range("I136").Value = 115000.445
$ws.Range("K136").Value = 345001.335
$ws.Range("M136").Value = -342451.335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 10000
$ws.Range("J31").Value = 10000
$ws.Range("L31").Value = 10000
$ws.Range("N31").Value = -10696
$ws.Range("H45").Value = 7569
$ws.Range("I45").Value = 7569
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 7569
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -7078
$ws.Range("N45").Value = ""
$ws.Range("H50").Value = 23084
$ws.Range("J50").Value = 23084
$ws.Range("L50").Value = 23084
$ws.Range("N50").Value = -24346
$ws.Range("H61").Value = 19994.75
$ws.Range("I61").Value = 19994.75
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 19994.75
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -19702.75
$ws.Range("N61").Value = ""
$ws.Range("H81").Value = 1237.6666
$ws.Range("J81").Value = 762.5
$ws.Range("L81").Value = 1525
$ws.Range("N81").Value = -3647
$ws.Range("H84").Value = 1237.6666
$ws.Range("J84").Value = 762.5
$ws.Range("L84").Value = 7625
$ws.Range("N84").Value = -18233
$ws.Range("H122").Value = 3094.3333
$ws.Range("I122").Value = 1454.75
$ws.Range("J122").Value = 6373.5
$ws.Range("K122").Value = 4364.25
$ws.Range("L122").Value = 19120.5
$ws.Range("M122").Value = -1914.25
$ws.Range("N122").Value = -24020.5
$ws.Range("H132").Value = 58408.277
$ws.Range("I132").Value = 9964.4
$ws.Range("J132").Value = 300627.66
$ws.Range("K132").Value = 29893.2
$ws.Range("L132").Value = 901882.98
$ws.Range("M132").Value = -27363.2
$ws.Range("N132").Value = -906942.98
$ws.Range("H136").Value = 4209.8423
$ws.Range("I136").Value = 3544.2727
$ws.Range("J136").Value = 5125
$ws.Range("K136").Value = 10632.8181
$ws.Range("L136").Value = 15375
$ws.Range("M136").Value = -8082.8181
$ws.Range("N136").Value = -20475
